# Update "Initial issue with label Testcase" config:
#  - rework Sheet1 (config) table: new header row (Project Name / URL / Note),
#    four project rows (erp-server, erp-web, xm-web, xm-api) with hyperlinked URLs
#  - add a second, empty "Sheet2"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row ---
$ws.Range("A1").Value = "Project Name"
$ws.Range("B1").Value = "URL"
$ws.Range("C1").Value = "Note"

# --- Data rows ---
$ws.Range("A2").Value = "erp-server"
$ws.Range("B2").Value = "https://git.iptp.net/erp/erp-server/-/issues "
$ws.Range("C2").Value = "Erp Server Project"

$ws.Range("A3").Value = "erp-web"
$ws.Range("B3").Value = "https://git.iptp.net/erp/erp-web/-/issues"
$ws.Range("C3").Value = "Erp Web Project"

$ws.Range("A4").Value = "xm-web"
$ws.Range("B4").Value = "https://git.iptp.net/xm/xm-web/-/issues"
$ws.Range("C4").Value = "XM Web Project"

$ws.Range("A5").Value = "xm-api"
$ws.Range("B5").Value = "https://git.iptp.net/andre/xm-api/-/issues"
$ws.Range("C5").Value = "XM API Project"

# --- Hyperlinks (added in this order so generated relationship ids line up) ---
$ws.Hyperlinks.Add($ws.Range("B3"), "https://git.iptp.net/erp/erp-web/-/issues") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://git.iptp.net/erp/erp-server/-/issues ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://git.iptp.net/xm/xm-web/-/issues") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://git.iptp.net/andre/xm-api/-/issues") | Out-Null

# --- Selection marker left by the author on Sheet1 ---
$ws.Range("E18").Select() | Out-Null

# --- Add a second, empty worksheet after Sheet1 ---
$sheet2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheet2.Name = "Sheet2"

# Keep Sheet1 as the active/selected sheet
$ws.Activate()
